$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# --- Shift the "Date" column (A2:A91) forward by one day -----------------
# Each row's date becomes the date that used to be one row below it, i.e.
# row N now shows what row N+1 used to show. Using Copy/Paste (rather than
# assigning the string via Value2) preserves the cells' text typing so
# Excel does not reinterpret the ISO-looking "yyyy-MM-dd" strings as date
# serial numbers.
$ws.Range("A3:A91").Copy($ws.Range("A2:A90"))

# Build the brand new trailing date ("2026-02-15") as literal text using a
# scratch cell + formula (so it is never passed through the Value2 setter,
# which auto-converts recognizable date strings into date serials). Then
# copy only the *value* into A91 and clean up the scratch cell so no extra
# styles/number formats linger in the workbook.
$scratch = $ws.Range("Z1")
$scratch.Formula = "=""2026-02-15"""
$scratch.Copy()
$ws.Range("A91").PasteSpecial(-4163)
$scratch.Clear()

# --- Shift the "HTTPS URLs" values (C2:C91) forward by one day -----------
$oldC = @()
for ($r = 2; $r -le 91; $r++) {
    $oldC += $ws.Cells.Item($r, 3).Value2
}
for ($r = 2; $r -le 90; $r++) {
    $ws.Cells.Item($r, 3).Value2 = $oldC[$r - 1]
}
$ws.Cells.Item(91, 3).Value2 = 31
